# Auto-generated Excel COM-interop script applying the Malboro_Profits value updates.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1765.1428
$ws.Range("J32").Value = 1841.4
$ws.Range("L32").Value = 1841.4
$ws.Range("N32").Value = -2493.4
$ws.Range("H33").Value = 19232802
$ws.Range("I33").Value = 22728040
$ws.Range("K33").Value = 22728040
$ws.Range("M33").Value = -22727811
$ws.Range("H41").Value = 562.8889
$ws.Range("I41").Value = 344.66666
$ws.Range("K41").Value = 344.66666
$ws.Range("M41").Value = 95.33334000000002
$ws.Range("H64").Value = 83339496
$ws.Range("J64").Value = 7499.75
$ws.Range("L64").Value = 7499.75
$ws.Range("N64").Value = -7995.75
$ws.Range("H67").Value = 83339496
$ws.Range("J67").Value = 7499.75
$ws.Range("L67").Value = 7499.75
$ws.Range("N67").Value = -9215.75
$ws.Range("H70").Value = 6208.222
$ws.Range("I70").Value = 6196.4287
$ws.Range("K70").Value = 18589.2861
$ws.Range("M70").Value = -18319.2861
$ws.Range("H73").Value = 6208.222
$ws.Range("I73").Value = 6196.4287
$ws.Range("K73").Value = 18589.2861
$ws.Range("M73").Value = -17653.2861
$ws.Range("H76").Value = 10905.5
$ws.Range("I76").Value = 3858.25
$ws.Range("K76").Value = 3858.25
$ws.Range("M76").Value = -3543.25
$ws.Range("H79").Value = 10905.5
$ws.Range("I79").Value = 3858.25
$ws.Range("K79").Value = 3858.25
$ws.Range("M79").Value = -2766.25
$ws.Range("H116").Value = 6153.8887
$ws.Range("I116").Value = 6235.625
$ws.Range("K116").Value = 6235.625
$ws.Range("M116").Value = -2793.625
$ws.Range("H132").Value = 22250.643
$ws.Range("I132").Value = 24838
$ws.Range("K132").Value = 74514
$ws.Range("M132").Value = -71984
$ws.Range("H137").Value = 8904.375
$ws.Range("I137").Value = 2916.1765
$ws.Range("K137").Value = 8748.529500000001
$ws.Range("M137").Value = -6198.529500000001
$ws.Range("H138").Value = 1828.59
$ws.Range("I138").Value = 1197.5
$ws.Range("J138").Value = 2074.014
$ws.Range("K138").Value = 3592.5
$ws.Range("L138").Value = 6222.042
$ws.Range("M138").Value = 1547.5
$ws.Range("N138").Value = -16502.042
$ws.Range("H141").Value = 4043.1428
$ws.Range("I141").Value = 4088.5625
$ws.Range("K141").Value = 12265.6875
$ws.Range("M141").Value = -7085.6875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3183.7693
$ws.Range("I45").Value = 3255.9473
$ws.Range("J45").Value = 2987.8572
$ws.Range("K45").Value = 3255.9473
$ws.Range("L45").Value = 2987.8572
$ws.Range("M45").Value = -2878.9473
$ws.Range("N45").Value = -3741.8572
$ws.Range("H74").Value = 12627.447
$ws.Range("I74").Value = 2117.7036
$ws.Range("J74").Value = 38424.09
$ws.Range("K74").Value = 2117.7036
$ws.Range("L74").Value = 38424.09
$ws.Range("M74").Value = -1243.7036
$ws.Range("N74").Value = -40172.09
$ws.Range("H77").Value = 12627.447
$ws.Range("I77").Value = 2117.7036
$ws.Range("J77").Value = 38424.09
$ws.Range("K77").Value = 10588.518
$ws.Range("L77").Value = 192120.45
$ws.Range("M77").Value = -6220.518
$ws.Range("N77").Value = -200856.45
$ws.Range("H88").Value = 1813.2941
$ws.Range("I88").Value = 1390.6666
$ws.Range("K88").Value = 1390.6666
$ws.Range("M88").Value = -984.6666
$ws.Range("H91").Value = 1813.2941
$ws.Range("I91").Value = 1390.6666
$ws.Range("K91").Value = 1390.6666
$ws.Range("M91").Value = 13.33339999999998
$ws.Range("H101").Value = 38999.5
$ws.Range("J101").Value = 38999.5
$ws.Range("L101").Value = 38999.5
$ws.Range("N101").Value = -45489.5
$ws.Range("H110").Value = 5323.923
$ws.Range("I110").Value = 6889.8887
$ws.Range("K110").Value = 6889.8887
$ws.Range("M110").Value = -4844.8887
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 29458.725
$ws.Range("I20").Value = 9944.75
$ws.Range("K20").Value = 9944.75
$ws.Range("M20").Value = -9697.75
$ws.Range("H94").Value = 1043.9344
$ws.Range("I94").Value = 554.8605
$ws.Range("K94").Value = 554.8605
$ws.Range("M94").Value = -103.8605
$ws.Range("H103").Value = 13442.8
$ws.Range("J103").Value = 11803.5
$ws.Range("L103").Value = 11803.5
$ws.Range("N103").Value = -14147.5
$ws.Range("H134").Value = 47928.703
$ws.Range("I134").Value = 53052.35
$ws.Range("J134").Value = 33289.715
$ws.Range("K134").Value = 159157.05
$ws.Range("L134").Value = 99869.14499999999
$ws.Range("M134").Value = -156622.05
$ws.Range("N134").Value = -104939.145
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19229.586
$ws.Range("I31").Value = 8075.9375
$ws.Range("K31").Value = 8075.9375
$ws.Range("M31").Value = -7780.9375
$ws.Range("H34").Value = 19229.586
$ws.Range("I34").Value = 8075.9375
$ws.Range("K34").Value = 8075.9375
$ws.Range("M34").Value = -7873.9375
$ws.Range("H134").Value = 7485.423
$ws.Range("I134").Value = 2472.625
$ws.Range("K134").Value = 7417.875
$ws.Range("M134").Value = -4882.875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 607.375
$ws.Range("I23").Value = 600
$ws.Range("K23").Value = 1800
$ws.Range("M23").Value = -1565
$ws.Range("H33").Value = 144.16667
$ws.Range("I33").Value = 97
$ws.Range("J33").Value = 191.33333
$ws.Range("K33").Value = 582
$ws.Range("L33").Value = 1147.99998
$ws.Range("M33").Value = -299
$ws.Range("N33").Value = -1713.99998
$ws.Range("H107").Value = 2119.5
$ws.Range("I107").Value = 567.4545000000001
$ws.Range("J107").Value = 2932.476
$ws.Range("K107").Value = 1702.3635
$ws.Range("L107").Value = 8797.428
$ws.Range("M107").Value = 217.6364999999998
$ws.Range("N107").Value = -12637.428
$ws.Range("H114").Value = 1122.25
$ws.Range("I114").Value = 99.25
$ws.Range("J114").Value = 2145.25
$ws.Range("K114").Value = 297.75
$ws.Range("L114").Value = 6435.75
$ws.Range("M114").Value = 2956.25
$ws.Range("N114").Value = -12943.75
$ws.Range("H117").Value = 608.55554
$ws.Range("I117").Value = 496.85715
$ws.Range("J117").Value = 999.5
$ws.Range("K117").Value = 1490.57145
$ws.Range("L117").Value = 2998.5
$ws.Range("M117").Value = 1951.42855
$ws.Range("N117").Value = -9882.5
$ws.Range("H131").Value = 1449.47
$ws.Range("I131").Value = 1050
$ws.Range("J131").Value = 1474.9681
$ws.Range("K131").Value = 3150
$ws.Range("L131").Value = 4424.9043
$ws.Range("M131").Value = 1890
$ws.Range("N131").Value = -14504.9043
$ws.Range("H140").Value = 3856.5715
$ws.Range("I140").Value = 3999.5
$ws.Range("J140").Value = 2999
$ws.Range("K140").Value = 11998.5
$ws.Range("L140").Value = 8997
$ws.Range("M140").Value = -6818.5
$ws.Range("N140").Value = -19357
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11063.71
$ws.Range("I80").Value = 11272.056
$ws.Range("K80").Value = 11272.056
$ws.Range("M80").Value = -10274.056
$ws.Range("H83").Value = 11063.71
$ws.Range("I83").Value = 11272.056
$ws.Range("K83").Value = 56360.28
$ws.Range("M83").Value = -51368.28
$ws.Range("H97").Value = 1684.6875
$ws.Range("I97").Value = 1728.2307
$ws.Range("K97").Value = 1728.2307
$ws.Range("M97").Value = -1232.2307
$ws.Range("H135").Value = 95109.5
$ws.Range("J135").Value = 95109.5
$ws.Range("L135").Value = 95109.5
$ws.Range("N135").Value = -105249.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2446.2163
$ws.Range("I22").Value = 1914.0869
$ws.Range("K22").Value = 1914.0869
$ws.Range("M22").Value = -1619.0869
$ws.Range("H27").Value = 2446.2163
$ws.Range("I27").Value = 1914.0869
$ws.Range("K27").Value = 1914.0869
$ws.Range("M27").Value = -1807.0869
$ws.Range("H46").Value = 3680.3845
$ws.Range("I46").Value = 2199
$ws.Range("J46").Value = 3803.8333
$ws.Range("K46").Value = 2199
$ws.Range("L46").Value = 3803.8333
$ws.Range("M46").Value = -2011
$ws.Range("N46").Value = -4179.8333
$ws.Range("H61").Value = 3597
$ws.Range("I61").Value = 2830.28
$ws.Range("K61").Value = 2830.28
$ws.Range("M61").Value = -2628.28
$ws.Range("H82").Value = 4419.6665
$ws.Range("I82").Value = 6630
$ws.Range("J82").Value = 3314.5
$ws.Range("K82").Value = 6630
$ws.Range("L82").Value = 3314.5
$ws.Range("M82").Value = -6269
$ws.Range("N82").Value = -4036.5
$ws.Range("H85").Value = 4419.6665
$ws.Range("I85").Value = 6630
$ws.Range("J85").Value = 3314.5
$ws.Range("K85").Value = 6630
$ws.Range("L85").Value = 3314.5
$ws.Range("M85").Value = -5382
$ws.Range("N85").Value = -5810.5
$ws.Range("H100").Value = 5598.2666
$ws.Range("I100").Value = 9479.200000000001
$ws.Range("J100").Value = 3657.8
$ws.Range("K100").Value = 9479.200000000001
$ws.Range("L100").Value = 3657.8
$ws.Range("M100").Value = -8938.200000000001
$ws.Range("N100").Value = -4739.8
$ws.Range("H113").Value = 3597
$ws.Range("I113").Value = 2830.28
$ws.Range("K113").Value = 2830.28
$ws.Range("M113").Value = -660.2800000000002
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
